# Auto commit at 2025-11-15 8:35:14.12
# Refresh the latest metrics pull: update the raw figures on "Metrics"
# (the "today" sheet's formulas pull from these and recalculate
# automatically), then restore the active-cell selection left on each
# sheet.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# Updated metric values (rows 2-13 of column B)
$metrics.Range("B2").Value  = 186917.99000000005
$metrics.Range("B3").Value  = 164240.72
$metrics.Range("B4").Value  = 58091.86
$metrics.Range("B5").Value  = 7697
$metrics.Range("B6").Value  = 4983163.7400000012
$metrics.Range("B7").Value  = 4206317.4000000004
$metrics.Range("B8").Value  = 1465051.69
$metrics.Range("B9").Value  = 193904
$metrics.Range("B10").Value = 33448544.730000008
$metrics.Range("B11").Value = 31481592.559999999
$metrics.Range("B12").Value = 11746773.73
$metrics.Range("B13").Value = 1291534

# Leave the selection on Metrics at E13 ...
$metrics.Range("E13").Select() | Out-Null

# ... and the selection on "today" (the tab that stays active) at E7.
$today.Activate() | Out-Null
$today.Range("E7").Select() | Out-Null
